$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new prompt text to cell F7
$ws.Range("F7").Value = "推荐书/Provide a list of 20 most related best books with intro, in this format:   {book name}/{intro}.  Final output are in the following format:     - item 1     - item 2     - item 3"

# Update the view: scroll so D1 is the top-left cell, and select F7
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("F7").Select()
